# Applies the "No more default assets" change:
# - Config is now the only way to add energy assets (for households & industry).
# - The shared string "OTHER_ELECTRICITY_CONSUMPTION" is replaced everywhere by
#   "ELECTRICITY_CONSUMPTION_PROFILE" on the consumptionAssets sheet.
# - A new "vehicle_scaling" (column P) value is recorded for every row on the
#   storageAssets sheet.
# - Various sheet selections / the active sheet are updated to reflect where the
#   author ended up after editing.

$wb = $excel.ActiveWorkbook

$wsConsumption = $wb.Worksheets.Item("consumptionAssets")
$wsProduction  = $wb.Worksheets.Item("productionAssets")
$wsConversion  = $wb.Worksheets.Item("conversionAssets")
$wsStorage     = $wb.Worksheets.Item("storageAssets")

# ---------------------------------------------------------------------------
# consumptionAssets: rename the "OTHER_ELECTRICITY_CONSUMPTION" energyAssetType
# to "ELECTRICITY_CONSUMPTION_PROFILE" for the rows that used it.
# ---------------------------------------------------------------------------
$wsConsumption.Range("E2").Value = "ELECTRICITY_CONSUMPTION_PROFILE"
$wsConsumption.Range("E4").Value = "ELECTRICITY_CONSUMPTION_PROFILE"
$wsConsumption.Range("E5").Value = "ELECTRICITY_CONSUMPTION_PROFILE"

# ---------------------------------------------------------------------------
# storageAssets: add the new "vehicle_scaling" column (P) values.
# ---------------------------------------------------------------------------
$wsStorage.Range("P2").Value = 0
$wsStorage.Range("P3").Value = 0
$wsStorage.Range("P4").Value = 0
$wsStorage.Range("P5").Value = 0
$wsStorage.Range("P6").Value = 0
$wsStorage.Range("P7").Value = 0
$wsStorage.Range("P8").Value = 0
$wsStorage.Range("P9").Value = 0
$wsStorage.Range("P10").Value = 1
$wsStorage.Range("P11").Value = 0
$wsStorage.Range("P12").Value = 0
$wsStorage.Range("P13").Value = 0
$wsStorage.Range("P14").Value = 0
$wsStorage.Range("P15").Value = 1
$wsStorage.Range("P16").Value = 0

# ---------------------------------------------------------------------------
# Sheet selections, matching where the author left the cursor on each tab.
# ---------------------------------------------------------------------------
$wsConsumption.Activate()
$wsConsumption.Range("E28").Select()

$wsProduction.Activate()
$wsProduction.Range("C25").Select()

$wsConversion.Activate()
$wsConversion.Range("K15").Select()

# storageAssets is the sheet the author finished on, so it is activated last
# and becomes the selected tab.
$wsStorage.Activate()
$wsStorage.Range("P16").Select()
